# Apply the crypto price / volume(1h) refresh described in the commit diff.
# Only cells B/C/D/E for the affected rows change; column A (rank index) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.110.10'
$ws.Range('E2').Value = '  +3.04%  '

# Row 3
$ws.Range('D3').Value = '3.416.07'
$ws.Range('E3').Value = '  +3.71%  '

# Row 4
$ws.Range('E4').Value = '  -0.20%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '577.92'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  +2.96%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '139.22'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +7.71%  '

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range('E7').Value = '  -0.14%  '

# Row 8
$ws.Range('D8').Value = '3.411.32'

# Row 9
$ws.Range('E9').Value = '  +1.46%  '

# Row 10
$ws.Range('E10').Value = '  +2.35%  '

# Row 11
$ws.Range('E11').Value = '  +9.98%  '

# Row 12
$ws.Range('E12').Value = '  +6.89%  '

# Row 13
$ws.Range('D13').Value = '3.991.24'
$ws.Range('E13').Value = '  +3.40%  '

# Row 14
$ws.Range('E14').Value = '  +2.02%  '

# Row 15
$ws.Range('E15').Value = '  +8.92%  '

# Row 16
$ws.Range('D16').Value = '3.415.29'
$ws.Range('E16').Value = '  +3.58%  '

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.59'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  +5.90%  '

# Row 18
$ws.Range('D18').Value = '62.084.40'
$ws.Range('E18').Value = '  +2.58%  '

# Row 19
$ws.Range('E19').Value = '  +6.62%  '

# Row 20
$ws.Range('E20').Value = '  +4.75%  '

# Row 21
$ws.Range('E21').Value = '  +6.42%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '392.30'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  +12.18%  '

# Row 23
$ws.Range('E23').Value = '  +4.05%  '

# Row 24
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000130'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  +20.53%  '

# Row 25
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.548.74'
$ws.Range('E25').Value = '  +3.49%  '

# Row 26
$ws.Range('E26').Value = '  +0.06%  '

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.80'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  +3.82%  '

# Row 28
$ws.Range('E28').Value = '  +10.72%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.66'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  +4.53%  '

# Row 30
$ws.Range('E30').Value = '  +0.02%  '

# Row 31
$ws.Range('E31').Value = '  +6.81%  '

# Row 32
$ws.Range('E32').Value = '  +5.56%  '

# Row 33
$ws.Range('E33').Value = '  +3.70%  '

# Row 34
$ws.Range('D34').Value = '3.440.81'
$ws.Range('E34').Value = '  +3.44%  '

# Row 35
$ws.Range('E35').Value = '  -0.04%  '

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.66'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  +4.87%  '

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.55'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  +4.36%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.00'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  +3.71%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.57'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +6.62%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '161.82'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +2.93%  '

# Row 42
$ws.Range('E42').Value = '  +14.70%  '

# Row 43
$ws.Range('B43').Value = 'ONDO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.24'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  +7.00%  '

# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -0.05%  '

# Row 45
$ws.Range('E45').Value = '  +5.31%  '

# Row 46
$ws.Range('E46').Value = '  +4.15%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.29'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  +11.80%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '41.62'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  +1.83%  '

# Row 49
$ws.Range('E49').Value = '  +5.09%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.03'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  +6.39%  '

# Row 51
$ws.Range('D51').Value = '2.400.77'
$ws.Range('E51').Value = '  +11.16%  '
